# Third Commit - Bruins Forwards Data
# The "Rk" column (column B) is removed from the goalies sheet; every
# column to its right shifts one position to the left (C->B, D->C, ...,
# AD->AC). Deleting the column itself performs the shift plus re-indexes
# the sheet's used range/dimension for us.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).Delete()
